$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.504.21'
$ws.Range("E2").Value = '  -2.30%  '
$ws.Range("D3").Value = '2.590.24'
$ws.Range("E3").Value = '  -3.17%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.16'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.53'
$ws.Range("E6").Value = '  -2.52%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.597'
$ws.Range("E8").Value = '  -2.13%  '
$ws.Range("D9").Value = '2.600.08'
$ws.Range("E9").Value = '  -2.75%  '
$ws.Range("E10").Value = '  -2.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("E12").Value = '  +10.27%  '
$ws.Range("E13").Value = '  +4.09%  '
$ws.Range("D14").Value = '3.048.32'
$ws.Range("E14").Value = '  -2.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.39'
$ws.Range("E15").Value = '  +6.70%  '
$ws.Range("D16").Value = '59.427.33'
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '2.595.66'
$ws.Range("E18").Value = '  -2.77%  '
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '339.06'
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.43'
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.52'
$ws.Range("E22").Value = '  +2.21%  '
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.70'
$ws.Range("E24").Value = '  -4.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.473'
$ws.Range("E25").Value = '  +7.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.162'
$ws.Range("E27").Value = '  -1.53%  '
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.19'
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("E32").Value = '  -2.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '158.48'
$ws.Range("E33").Value = '  +1.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.11'
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.18'
$ws.Range("E36").Value = '  +0.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.895'
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.873'
$ws.Range("E38").Value = '  -3.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.41'
$ws.Range("E39").Value = '  -0.56%  '
$ws.Range("E40").Value = '  -2.00%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '295.27'
$ws.Range("E41").Value = '  -2.46%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.69'
$ws.Range("E42").Value = '  +0.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '139.68'
$ws.Range("E43").Value = '  +9.11%  '
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0977'
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0533'
$ws.Range("E49").Value = '  -2.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.83'
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("D51").Value = '1.964.51'
$ws.Range("E51").Value = '  -0.20%  '
